$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.2456140350877193
$ws.Range("C2").Value = 0.4502923976608187
$ws.Range("J2").Value = 0.01169590643274854
$ws.Range("P2").Value = 0.1286549707602339
$ws.Range("S2").Value = 0.1637426900584795
$ws.Range("B3").Value = 0.01282051282051282
$ws.Range("C3").Value = 0.01282051282051282
$ws.Range("J3").Value = 0.02564102564102564
$ws.Range("P3").Value = 0.6410256410256411
$ws.Range("S3").Value = 0.3076923076923077
$ws.Range("P4").Value = 0.7307692307692307
$ws.Range("S4").Value = 0.2692307692307692
$ws.Range("B6").Value = 0.05555555555555555
$ws.Range("D6").Value = 0.01111111111111111
$ws.Range("F6").Value = 0.03333333333333333
$ws.Range("J6").Value = 0.2555555555555555
$ws.Range("O6").Value = 0.01111111111111111
$ws.Range("Q6").Value = 0.1444444444444444
$ws.Range("R6").Value = 0.1
$ws.Range("S6").Value = 0.3888888888888889
$ws.Range("B7").Value = 0.1111111111111111
$ws.Range("D7").Value = 0.01587301587301587
$ws.Range("F7").Value = 0.07936507936507936
$ws.Range("J7").Value = 0.07936507936507936
$ws.Range("O7").Value = 0.01587301587301587
$ws.Range("Q7").Value = 0.1587301587301587
$ws.Range("R7").Value = 0.06349206349206349
$ws.Range("S7").Value = 0.4761904761904762
$ws.Range("B8").Value = 0.07373271889400922
$ws.Range("D8").Value = 0.0184331797235023
$ws.Range("F8").Value = 0.05069124423963134
$ws.Range("J8").Value = 0.08294930875576037
$ws.Range("Q8").Value = 0.1244239631336406
$ws.Range("R8").Value = 0.1244239631336406
$ws.Range("S8").Value = 0.5253456221198156
$ws.Range("B9").Value = 0.09278350515463918
$ws.Range("F9").Value = 0.05154639175257732
$ws.Range("J9").Value = 0.06185567010309279
$ws.Range("O9").Value = 0.02061855670103093
$ws.Range("Q9").Value = 0.1649484536082474
$ws.Range("R9").Value = 0.08247422680412371
$ws.Range("S9").Value = 0.5257731958762887
$ws.Range("B10").Value = 0.1444444444444444
$ws.Range("D10").Value = 0.03333333333333333
$ws.Range("F10").Value = 0.05555555555555555
$ws.Range("J10").Value = 0.09682539682539683
$ws.Range("O10").Value = 0.02063492063492063
$ws.Range("Q10").Value = 0.2
$ws.Range("R10").Value = 0.05396825396825397
$ws.Range("S10").Value = 0.3952380952380952
$ws.Range("G11").Value = 0.11
$ws.Range("J11").Value = 0.12
$ws.Range("K11").Value = 0.17
$ws.Range("L11").Value = 0.59
$ws.Range("S11").Value = 0.01
$ws.Range("G12").Value = 0.7142857142857143
$ws.Range("J12").Value = 0.2063492063492063
$ws.Range("L12").Value = 0.03174603174603174
$ws.Range("S12").Value = 0.04761904761904762
$ws.Range("G13").Value = 0.4375
$ws.Range("J13").Value = 0.375
$ws.Range("S13").Value = 0.1875
$ws.Range("F15").Value = 0.009615384615384616
$ws.Range("H15").Value = 0.1538461538461539
$ws.Range("I15").Value = 0.0576923076923077
$ws.Range("J15").Value = 0.4038461538461539
$ws.Range("K15").Value = 0.07692307692307693
$ws.Range("O15").Value = 0.04807692307692308
$ws.Range("F16").Value = 0.06593406593406594
$ws.Range("H16").Value = 0.1318681318681319
$ws.Range("I16").Value = 0.05494505494505494
$ws.Range("J16").Value = 0.4615384615384616
$ws.Range("K16").Value = 0.0989010989010989
$ws.Range("O16").Value = 0.06593406593406594
$ws.Range("S16").Value = 0.1208791208791209
$ws.Range("F17").Value = 0.02105263157894737
$ws.Range("H17").Value = 0.1631578947368421
$ws.Range("I17").Value = 0.1526315789473684
$ws.Range("J17").Value = 0.3684210526315789
$ws.Range("K17").Value = 0.07894736842105263
$ws.Range("M17").Value = 0.01578947368421053
$ws.Range("O17").Value = 0.06315789473684211
$ws.Range("S17").Value = 0.1368421052631579
$ws.Range("F18").Value = 0.01219512195121951
$ws.Range("H18").Value = 0.1707317073170732
$ws.Range("I18").Value = 0.04878048780487805
$ws.Range("J18").Value = 0.475609756097561
$ws.Range("K18").Value = 0.08536585365853659
$ws.Range("M18").Value = 0.01219512195121951
$ws.Range("O18").Value = 0.03658536585365853
$ws.Range("S18").Value = 0.1585365853658537
$ws.Range("F19").Value = 0.0100864553314121
$ws.Range("H19").Value = 0.2118155619596542
$ws.Range("I19").Value = 0.07780979827089338
$ws.Range("J19").Value = 0.4265129682997118
$ws.Range("K19").Value = 0.06195965417867435
$ws.Range("M19").Value = 0.01873198847262248
$ws.Range("N19").Value = 0.001440922190201729
$ws.Range("O19").Value = 0.06772334293948126
$ws.Range("S19").Value = 0.1239193083573487
